$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (16) to the faculties table with a new "Others" entry
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Others"
$ws.Range("C16").Value = "Others"
